$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 (Ví điện tử account)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("A3").Value = "abcde"
$ws.Range("B3").Value = "122300787"

# Add new row 4 for the bank card account
$ws.Range("B4").NumberFormat = "@"
$ws.Range("A4").Value = "zcx"
$ws.Range("B4").Value = "2423"
$ws.Range("C4").Value = "Thẻ ngân hàng"
